# Adds season-record columns (Wins, Losses, Ties) to the roster sheet,
# mirroring the data scraped separately for the team's overall record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------
# Copy the formatting of the last existing header cell (bold, bordered,
# centered) onto the three new header cells, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows --------------------------------------------------------
# Every player row (2 through 52) gets the team's season record.
$wins = 66
$losses = 95
$ties = 0

for ($row = 2; $row -le 52; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins    # AD
    $ws.Cells.Item($row, 31).Value = $losses  # AE
    $ws.Cells.Item($row, 32).Value = $ties    # AF
}
